$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138 (shifts rows 138..236 down to 139..237)
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row with the new record
$ws.Cells.Item(138, 1).Value  = 6
$ws.Cells.Item(138, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(138, 3).Value  = "Metropolitana"
$ws.Cells.Item(138, 4).Value  = 44777
$ws.Cells.Item(138, 5).Value  = 13
$ws.Cells.Item(138, 6).Value  = 100112022
$ws.Cells.Item(138, 7).Value  = "Arveja Verde"
$ws.Cells.Item(138, 8).Value  = "Perfection"
$ws.Cells.Item(138, 9).Value  = "Primera"
$ws.Cells.Item(138, 10).Value = 400
$ws.Cells.Item(138, 11).Value = 38000
$ws.Cells.Item(138, 12).Value = 39000
$ws.Cells.Item(138, 13).Value = 38425
$ws.Cells.Item(138, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(138, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(138, 16).Value = 1537
$ws.Cells.Item(138, 17).Value = 25
$ws.Cells.Item(138, 18).Value = "Hortaliza"
